$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell values (A1, B1, C1)
$ws.Range("A1").Value = 156.82324786752349
$ws.Range("B1").Value = 5.1973514400349146
$ws.Range("C1").Value = 1.0876115459882583

# Update column widths (B and C)
$ws.Columns.Item(2).ColumnWidth = 10.7109375
$ws.Columns.Item(3).ColumnWidth = 11.7109375
